$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values were regenerated (commit: "regen save_data to use K
# instead of Strike#, regen std/mean, calc and write s_vals").
# Update G2:G8 with the new K values from the diff.
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 2
